$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Data")

# --- PayGrades table (rows 35-40): mark Grade C / Class rows as having a
# currency code now, and flip the "Runmode" (column G) flag to Y for the
# rows that got exercised. Order mirrors how the author likely keyed the
# cells in (matches the shared-string append order in the saved file).

# Class 3 / Class 4 rows now run (G35, G36 -> Y)
$ws.Range("G35").Value = "Y"
$ws.Range("G36").Value = "Y"

# Grade C row (40): fill in Currency + Maximum Salary first
$ws.Range("C40").Value = "LKR"
$ws.Range("E40").Value = "'five lakhs"

# Grade B row (38): fix the currency code + minimum salary, mark it run
$ws.Range("C38").Value = "NZD"
$ws.Range("D38").Value = "'350000.25"
$ws.Range("G38").Value = "Y"

# Row 39 also marked as run
$ws.Range("G39").Value = "Y"

# Back to Grade C row (40): fill in the Minimum Salary
$ws.Range("D40").Value = "'325625"

# Column B ("Grade Name") widened to fit the new content
$ws.Columns.Item(2).ColumnWidth = 19.5
